$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 45, shifting existing rows 45:183 down to 46:184
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the latest notification entry
$ws.Cells.Item(45, 18).Value = "balance your axis"
$ws.Cells.Item(45, 19).Value = "2024-09-21 10:34:04"
